$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update task names for rows 24-28 (column B)
$ws.Range("B24").Value = "Fix Save button"
$ws.Range("B25").Value = "Find a way to merge all html forms, entity, repository and service"
$ws.Range("B26").Value = "Fix ajax for terminal_ID"
$ws.Range("B27").Value = "Fix tabEdit.html page, merge all tabs in one html page"
$ws.Range("B28").Value = "Fix tabEntity for parsing of values"

# Update hour entries for these rows (columns AE-AH correspond to Apr 27-30)
$ws.Range("AE24").Value = 5
$ws.Range("AE25").Value = 3
$ws.Range("AF25").Value = 3
$ws.Range("AF26").Value = 5
$ws.Range("AG27").Value = 8
$ws.Range("AH28").Value = 8

# Update selected cell in the sheet view
$ws.Range("B28:D28").Select()

$wb.Save()
